# Commit: Add data for 2022-12-17
# Updates the "through" date on the report (Dec 08 -> Dec 09) and refreshes
# the carjacking counts for a number of neighborhood/month cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet tab to reflect the new "through" date.
$ws.Name = "Through 2022-12-09"

# Update the report title text in column B row 1 (shared string referenced there).
$ws.Range("B1").Value = "December 2022 (through December 09)"

# --- Updated / newly populated monthly counts per neighborhood row ---

# Row 2 - Garfield Park
$ws.Range("N2").Value = 1
$ws.Range("AX2").Value = 2
$ws.Range("BJ2").Value = 3

# Row 3 - Humboldt Park
$ws.Range("AX3").Value = 3
$ws.Range("BV3").Value = 2

# Row 5 - (data row 5)
$ws.Range("BV5").Value = 4

# Row 7 - Washington Heights
$ws.Range("B7").Value = 1
$ws.Range("BJ7").Value = 1

# Row 11
$ws.Range("BV11").Value = 1

# Row 14
$ws.Range("AL14").Value = 2
$ws.Range("AX14").Value = 3

# Row 15
$ws.Range("N15").Value = 5

# Row 16
$ws.Range("Z16").Value = 2

# Row 20
$ws.Range("BJ20").Value = 2

# Row 22
$ws.Range("N22").Value = 3

# Row 27
$ws.Range("AL27").Value = 1

# Row 28
$ws.Range("N28").Value = 1

# Row 40
$ws.Range("N40").Value = 1

# Row 41
$ws.Range("Z41").Value = 3

# Row 42
$ws.Range("BJ42").Value = 1

# Row 45 - Bridgeport
$ws.Range("B45").Value = 1

# Row 49
$ws.Range("BJ49").Value = 1

# Row 64
$ws.Range("N64").Value = 5

# Row 74
$ws.Range("Z74").Value = 1

# Row 88 - Old Town
$ws.Range("B88").Value = 2

# Row 90 - Pullman
$ws.Range("AX90").Value = 1
